# "implemented add owners functionality"
# Adds an "owner" (gmail) column with mailto hyperlinks, renames the
# Vehicle Number column, refreshes Brand/Model/FRV code data, and appends
# year / isAc / km / date columns to the cars sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Column A: "name" -> "gmail", values become mailto: hyperlinked emails.
# Hyperlinks are added in this exact row order (2,3,4,5,6,8,7) so the
# relationship ids line up with how the sheet was actually authored.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "gmail"

$owners = @{
    2 = "test1@gmail.com"
    3 = "test2@gmail.com"
    4 = "test1@gmail.com"
    5 = "test3@gmail.com"
    6 = "test1@gmail.com"
    8 = "test3@gmail.com"
    7 = "test2@gmail.com"
}
foreach ($r in @(2,3,4,5,6,8,7)) {
    $email = $owners[$r]
    $ws.Range("A$r").Value = $email
    $ws.Hyperlinks.Add($ws.Range("A$r"), "mailto:$email") | Out-Null
}

# ---------------------------------------------------------------------
# Column D header rename (data itself is unchanged).
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Vehicle Registration Number"

# ---------------------------------------------------------------------
# Brand / Model (B, C) and FRV code (E) refresh per row.
# ---------------------------------------------------------------------
$data = @{
    2 = @{ B = "Land Rover"; C = "Defender";    E = "ASK-01" }
    3 = @{ B = "TATA";       C = "Harrier";     E = "ASK-02" }
    4 = @{ B = "Mahindra";   C = "XUV 700";     E = "ASK-01" }
    5 = @{ B = "TATA";       C = "Nexon";       E = "ASK-12" }
    6 = @{ B = "Mahindra";   C = "Bolero";      E = "BSK-01" }
    7 = @{ B = "Land Rover"; C = "Defender";    E = "ASK-04" }
    8 = @{ B = "Land Rover"; C = "Range Rover"; E = "BSK-02" }
}
foreach ($r in 2..8) {
    $row = $data[$r]
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 5).Value = $row.E
}

# ---------------------------------------------------------------------
# New columns: I=year, J=isAc, K=km, L=date
# ---------------------------------------------------------------------
$ws.Range("I1").Value = "year"
$ws.Range("J1").Value = "isAc"
$ws.Range("K1").Value = "km"
$ws.Range("L1").Value = "date"

$rows = @{
    2 = @{ I = 2020; J = $true;  K = 300 }
    3 = @{ I = 2021; J = $true;  K = 3000 }
    4 = @{ I = 2022; J = $true;  K = 4500 }
    5 = @{ I = 2021; J = $true;  K = 2300 }
    6 = @{ I = 2020; J = $false; K = 6755 }
    7 = @{ I = 2019; J = $false; K = 4356 }
    8 = @{ I = 2023; J = $true;  K = 24674 }
}
foreach ($r in 2..8) {
    $row = $rows[$r]
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
}

# Give the date column its number format on the first cell only (so it
# maps onto the built-in short-date format, numFmtId 14, rather than
# minting a bespoke custom format), then fan that single style out to
# the rest of the column via copy/paste-special so every date cell
# shares one style record instead of getting its own.
$ws.Range("L2").NumberFormat = "mm-dd-yy"
$ws.Range("L2").Value = "3/19/2024"
$ws.Range("L2").Copy() | Out-Null
$ws.Range("L3:L8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
foreach ($r in 3..8) {
    $ws.Cells.Item($r, 12).Value = "3/19/2024"
}

$ws.Columns.Item(12).ColumnWidth = 10.0

# Restore the original active-cell look (top-left scrolled to column B,
# selection resting on B8) that the authoring session ended on.
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B8").Select() | Out-Null
